$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1917.4286
$ws.Range("I29").Value = 951.5
$ws.Range("J29").Value = 2303.8
$ws.Range("K29").Value = 2854.5
$ws.Range("L29").Value = 6911.400000000001
$ws.Range("M29").Value = -2573.5
$ws.Range("N29").Value = -7473.400000000001

$ws.Range("H33").Value = 327.91666
$ws.Range("I33").Value = 346.3158
$ws.Range("J33").Value = 258
$ws.Range("K33").Value = 346.3158
$ws.Range("L33").Value = 258
$ws.Range("M33").Value = -117.3158
$ws.Range("N33").Value = -716

$ws.Range("H38").Value = 917.55554
$ws.Range("I38").Value = 149.66667
$ws.Range("J38").Value = 2453.3333
$ws.Range("K38").Value = 449.00001
$ws.Range("L38").Value = 7359.999899999999
$ws.Range("M38").Value = -77.00001000000003
$ws.Range("N38").Value = -8103.999899999999

$ws.Range("H41").Value = 18519846
$ws.Range("I41").Value = 33334468
$ws.Range("J41").Value = 1570
$ws.Range("K41").Value = 33334468
$ws.Range("L41").Value = 1570
$ws.Range("M41").Value = -33334028
$ws.Range("N41").Value = -2450

$ws.Range("H58").Value = 3877.7222
$ws.Range("J58").Value = 5549.9165
$ws.Range("L58").Value = 16649.7495
$ws.Range("N58").Value = -16949.7495

$ws.Range("H64").Value = 3554.8293
$ws.Range("I64").Value = 3563.7058
$ws.Range("J64").Value = 3548.5417
$ws.Range("K64").Value = 3563.7058
$ws.Range("L64").Value = 3548.5417
$ws.Range("M64").Value = -3315.7058
$ws.Range("N64").Value = -4044.5417

$ws.Range("H67").Value = 3554.8293
$ws.Range("I67").Value = 3563.7058
$ws.Range("J67").Value = 3548.5417
$ws.Range("K67").Value = 3563.7058
$ws.Range("L67").Value = 3548.5417
$ws.Range("M67").Value = -2705.7058
$ws.Range("N67").Value = -5264.5417

$ws.Range("H98").Value = 2383.0557
$ws.Range("I98").Value = 2055
$ws.Range("J98").Value = 3236
$ws.Range("K98").Value = 2055
$ws.Range("L98").Value = 3236
$ws.Range("M98").Value = -557
$ws.Range("N98").Value = -6232

$ws.Range("H112").Value = 2073.7173
$ws.Range("J112").Value = 2117.9773
$ws.Range("L112").Value = 6353.9319
$ws.Range("N112").Value = -8569.9319

$ws.Range("H113").Value = 2007.25
$ws.Range("I113").Value = 1722.5714
$ws.Range("K113").Value = 1722.5714
$ws.Range("M113").Value = 1531.4286

$ws.Range("H116").Value = 2442.0476
$ws.Range("I116").Value = 2104.8823
$ws.Range("J116").Value = 3875
$ws.Range("K116").Value = 2104.8823
$ws.Range("L116").Value = 3875
$ws.Range("M116").Value = 1337.1177
$ws.Range("N116").Value = -10759

$ws.Range("H122").Value = 2383.0557
$ws.Range("I122").Value = 2055
$ws.Range("J122").Value = 3236
$ws.Range("K122").Value = 6165
$ws.Range("L122").Value = 9708
$ws.Range("M122").Value = -3715
$ws.Range("N122").Value = -14608

$ws.Range("H129").Value = 884.0417
$ws.Range("I129").Value = 459.66666
$ws.Range("J129").Value = 944.6667
$ws.Range("K129").Value = 1378.99998
$ws.Range("L129").Value = 2834.0001
$ws.Range("M129").Value = 3621.00002
$ws.Range("N129").Value = -12834.0001

$ws.Range("H132").Value = 7581999.5
$ws.Range("I132").Value = 12348509
$ws.Range("J132").Value = 11660.235
$ws.Range("K132").Value = 37045527
$ws.Range("L132").Value = 34980.705
$ws.Range("M132").Value = -37042997
$ws.Range("N132").Value = -40040.705

$ws.Range("H138").Value = 726501.1
$ws.Range("J138").Value = 968248.8
$ws.Range("L138").Value = 2904746.4
$ws.Range("N138").Value = -2915026.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4430.2188
$ws.Range("I32").Value = 4581.8276
$ws.Range("J32").Value = 2964.6667
$ws.Range("K32").Value = 4581.8276
$ws.Range("L32").Value = 2964.6667
$ws.Range("M32").Value = -4294.8276
$ws.Range("N32").Value = -3538.6667

$ws.Range("H45").Value = 1472.625
$ws.Range("I45").Value = 1457.4667
$ws.Range("J45").Value = 1700
$ws.Range("K45").Value = 1457.4667
$ws.Range("L45").Value = 1700
$ws.Range("M45").Value = -1080.4667
$ws.Range("N45").Value = -2454

$ws.Range("H74").Value = 1092.6
$ws.Range("I74").Value = 995.0909
$ws.Range("K74").Value = 995.0909
$ws.Range("M74").Value = -121.0909

$ws.Range("H77").Value = 1092.6
$ws.Range("I77").Value = 995.0909
$ws.Range("K77").Value = 4975.4545
$ws.Range("M77").Value = -607.4544999999998

$ws.Range("H92").Value = 27887
$ws.Range("J92").Value = 27887
$ws.Range("L92").Value = 27887
$ws.Range("N92").Value = -32879

$ws.Range("H122").Value = 1461.091
$ws.Range("I122").Value = 1472
$ws.Range("J122").Value = 1448
$ws.Range("K122").Value = 4416
$ws.Range("L122").Value = 4344
$ws.Range("M122").Value = -1966
$ws.Range("N122").Value = -9244

$ws.Range("H132").Value = 2631.2917
$ws.Range("I132").Value = 2364.125
$ws.Range("J132").Value = 3165.625
$ws.Range("K132").Value = 7092.375
$ws.Range("L132").Value = 9496.875
$ws.Range("M132").Value = -4562.375
$ws.Range("N132").Value = -14556.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1687.1904
$ws.Range("I107").Value = 1249
$ws.Range("J107").Value = 2271.4443
$ws.Range("K107").Value = 1249
$ws.Range("L107").Value = 2271.4443
$ws.Range("M107").Value = 671
$ws.Range("N107").Value = -6111.4443

$ws.Range("H134").Value = 5399.16
$ws.Range("I134").Value = 1177.5555
$ws.Range("J134").Value = 16254.714
$ws.Range("K134").Value = 3532.6665
$ws.Range("L134").Value = 48764.142
$ws.Range("M134").Value = -997.6664999999998
$ws.Range("N134").Value = -53834.142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1535.45
$ws.Range("I31").Value = 1475.5
$ws.Range("J31").Value = 2075
$ws.Range("K31").Value = 1475.5
$ws.Range("L31").Value = 2075
$ws.Range("M31").Value = -1180.5
$ws.Range("N31").Value = -2665

$ws.Range("H34").Value = 1535.45
$ws.Range("I34").Value = 1475.5
$ws.Range("J34").Value = 2075
$ws.Range("K34").Value = 1475.5
$ws.Range("L34").Value = 2075
$ws.Range("M34").Value = -1273.5
$ws.Range("N34").Value = -2479

$ws.Range("H58").Value = 862.9231
$ws.Range("I58").Value = 850.3333
$ws.Range("J58").Value = 1014
$ws.Range("K58").Value = 850.3333
$ws.Range("L58").Value = 1014
$ws.Range("M58").Value = -647.3333
$ws.Range("N58").Value = -1420

$ws.Range("H94").Value = 1114.125
$ws.Range("I94").Value = 1341.3334
$ws.Range("J94").Value = 977.8
$ws.Range("K94").Value = 1341.3334
$ws.Range("L94").Value = 977.8
$ws.Range("M94").Value = -890.3334
$ws.Range("N94").Value = -1879.8

$ws.Range("H107").Value = 687.6
$ws.Range("I107").Value = 402.66666
$ws.Range("J107").Value = 877.55554
$ws.Range("K107").Value = 402.66666
$ws.Range("L107").Value = 877.55554
$ws.Range("M107").Value = 1517.33334
$ws.Range("N107").Value = -4717.55554

$ws.Range("H136").Value = 862.9231
$ws.Range("I136").Value = 850.3333
$ws.Range("J136").Value = 1014
$ws.Range("K136").Value = 2550.9999
$ws.Range("L136").Value = 3042
$ws.Range("M136").Value = -0.9998999999997977
$ws.Range("N136").Value = -8142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 18869278
$ws.Range("J131").Value = 1554.4546
$ws.Range("L131").Value = 4663.3638
$ws.Range("N131").Value = -14743.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1177.4166
$ws.Range("I122").Value = 1041.25
$ws.Range("J122").Value = 1449.75
$ws.Range("K122").Value = 3123.75
$ws.Range("L122").Value = 4349.25
$ws.Range("M122").Value = -673.75
$ws.Range("N122").Value = -9249.25

$ws.Range("H126").Value = 2961.45
$ws.Range("I126").Value = 1668.8
$ws.Range("J126").Value = 4254.1
$ws.Range("K126").Value = 5006.4
$ws.Range("L126").Value = 12762.3
$ws.Range("M126").Value = -2536.4
$ws.Range("N126").Value = -17702.3

$ws.Range("H132").Value = 2846.8572
$ws.Range("I132").Value = 3127.875
$ws.Range("J132").Value = 2472.1667
$ws.Range("K132").Value = 9383.625
$ws.Range("L132").Value = 7416.500100000001
$ws.Range("M132").Value = -6853.625
$ws.Range("N132").Value = -12476.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1015.2381
$ws.Range("I22").Value = 982.6667
$ws.Range("J22").Value = 1096.6666
$ws.Range("K22").Value = 982.6667
$ws.Range("L22").Value = 1096.6666
$ws.Range("M22").Value = -687.6667
$ws.Range("N22").Value = -1686.6666

$ws.Range("H27").Value = 1015.2381
$ws.Range("I27").Value = 982.6667
$ws.Range("J27").Value = 1096.6666
$ws.Range("K27").Value = 982.6667
$ws.Range("L27").Value = 1096.6666
$ws.Range("M27").Value = -875.6667
$ws.Range("N27").Value = -1310.6666

$ws.Range("H68").Value = 1329.1333
$ws.Range("I68").Value = 1139.1
$ws.Range("K68").Value = 1139.1
$ws.Range("M68").Value = -390.0999999999999

$ws.Range("H71").Value = 1329.1333
$ws.Range("I71").Value = 1139.1
$ws.Range("K71").Value = 5695.5
$ws.Range("M71").Value = -1951.5

$ws.Range("H93").Value = 689.3
$ws.Range("I93").Value = 591.8570999999999
$ws.Range("J93").Value = 916.6667
$ws.Range("K93").Value = 591.8570999999999
$ws.Range("L93").Value = 916.6667
$ws.Range("M93").Value = 656.1429000000001
$ws.Range("N93").Value = -3412.6667

$ws.Range("H100").Value = 1950
$ws.Range("I100").Value = 1566.6666
$ws.Range("K100").Value = 1566.6666
$ws.Range("M100").Value = -1025.6666

$ws.Range("H132").Value = 129236.375
$ws.Range("I132").Value = 6444.5
$ws.Range("J132").Value = 170167
$ws.Range("K132").Value = 19333.5
$ws.Range("L132").Value = 510501
$ws.Range("M132").Value = -16803.5
$ws.Range("N132").Value = -515561

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 7432381.5
$ws.Range("I122").Value = 8969484
$ws.Range("J122").Value = 3052.5
$ws.Range("K122").Value = 26908452
$ws.Range("L122").Value = 9157.5
$ws.Range("M122").Value = -26906002
$ws.Range("N122").Value = -14057.5

$ws.Range("H136").Value = 745
$ws.Range("I136").Value = 317.5
$ws.Range("K136").Value = 952.5
$ws.Range("M136").Value = 1597.5
